$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 188-189, shifting existing rows 188-197 down to 190-199.
$ws.Rows("188:189").Insert()

# New row 188: Santina / Especial, dated 45267 (2023-12-07)
$ws.Range("A188").Value = 7
$ws.Range("B188").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C188").Value = "Ñuble"
$ws.Range("D188").Value = 45267
$ws.Range("E188").Value = 16
$ws.Range("F188").Value = "Fruta"
$ws.Range("G188").Value = 100103
$ws.Range("H188").Value = "Frutos de hueso (carozo)"
$ws.Range("I188").Value = 100103001
$ws.Range("J188").Value = "Cereza"
$ws.Range("K188").Value = "Santina"
$ws.Range("L188").Value = "Especial"
$ws.Range("M188").Value = 150
$ws.Range("N188").Value = 14000
$ws.Range("O188").Value = 14000
$ws.Range("P188").Value = 14000
$ws.Range("Q188").Value = "$/bandeja 10 kilos"
$ws.Range("R188").Value = "Provincia de Curicó"
$ws.Range("S188").Value = 1400
$ws.Range("T188").Value = 10

# New row 189: Santina / Primera, dated 45267 (2023-12-07)
$ws.Range("A189").Value = 7
$ws.Range("B189").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C189").Value = "Ñuble"
$ws.Range("D189").Value = 45267
$ws.Range("E189").Value = 16
$ws.Range("F189").Value = "Fruta"
$ws.Range("G189").Value = 100103
$ws.Range("H189").Value = "Frutos de hueso (carozo)"
$ws.Range("I189").Value = 100103001
$ws.Range("J189").Value = "Cereza"
$ws.Range("K189").Value = "Santina"
$ws.Range("L189").Value = "Primera"
$ws.Range("M189").Value = 150
$ws.Range("N189").Value = 12000
$ws.Range("O189").Value = 12000
$ws.Range("P189").Value = 12000
$ws.Range("Q189").Value = "$/bandeja 10 kilos"
$ws.Range("R189").Value = "Provincia de Curicó"
$ws.Range("S189").Value = 1200
$ws.Range("T189").Value = 10
